$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# TEST_CASES sheet: fill in TC_OWNER_PATH (D) / TC_OWNER_ID (E) sample values
# ---------------------------------------------------------------------------
$wsTestCases = $wb.Worksheets.Item("TEST_CASES")

# Values are written in this precise order so that the new shared-string
# table entries are created in the same sequence as the target workbook.
$wsTestCases.Range("D2").Value = "path/row1"
$wsTestCases.Range("E2").Value = 11

$wsTestCases.Range("D3").Value = "path/row2"
$wsTestCases.Range("E3").Value = 12

$wsTestCases.Range("D5").Value = "path/row4"
$wsTestCases.Range("E5").Value = 14

$wsTestCases.Range("D6").Value = "path/row5"
$wsTestCases.Range("E6").Value = 15

$wsTestCases.Range("D7").Value = "path/row6"
$wsTestCases.Range("E7").Value = 16

$wsTestCases.Range("D8").Value = "path/row7"
$wsTestCases.Range("E8").Value = 17

$wsTestCases.Range("D9").Value = "path/row8"
$wsTestCases.Range("E9").Value = 18

$wsTestCases.Range("D10").Value = "path/row9"
$wsTestCases.Range("E10").Value = 19

$wsTestCases.Range("D11").Value = "path/row10"
$wsTestCases.Range("E11").Value = 20

$wsTestCases.Range("D4").Value = "path/row3"
$wsTestCases.Range("E4").Value = 13

# Re-apply the row-2 cell format (thin border style) onto D3:E11 so the newly
# filled cells match the look of the header/first data row.
$wsTestCases.Range("D2:E2").Copy()
$wsTestCases.Range("D3:E11").PasteSpecial(-4122)

# The worksheet no longer needs the trailing blank rows 12-16.
$wsTestCases.Rows("12:16").Delete()

# Widen column D slightly to fit the new "path/rowX" values.
$wsTestCases.Columns("D").ColumnWidth = 9.5

# ---------------------------------------------------------------------------
# STEPS sheet: move the selected cell
# ---------------------------------------------------------------------------
$wsSteps = $wb.Worksheets.Item("STEPS")
$wsSteps.Activate() | Out-Null
$wsSteps.Range("B34").Select() | Out-Null

# ---------------------------------------------------------------------------
# TEST_CASES sheet becomes the active tab again, with a new selection
# ---------------------------------------------------------------------------
$wsTestCases.Activate() | Out-Null
$wsTestCases.Range("C18").Select() | Out-Null
